$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B10").Value = "hlthdist_fctb_clst"
$ws.Range("B10").Select()
